$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing merged cells so we can freely rewrite the layout
$ws.Cells.UnMerge() | Out-Null

# Stash the header row's existing format (fontId=1, border, centered) well off
# to the side (column Z, untouched by the later column insert) before wiping
# everything, so we don't invent any new style combos later.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Wipe all contents + formatting from the previously used range (this drops the
# colored/merged "Flight" banner rows' styles entirely).
$ws.Range("A1:C7").Clear() | Out-Null

# Insert a new column before column B for "Squadron Assignment"
$ws.Columns.Item(2).Insert() | Out-Null

# Restore the header format onto the (now 4-wide) header row, then tidy up the
# scratch cell used to stash it. The earlier column insert shifted the stash
# from Z1 to AA1.
$ws.Range("AA1").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("AA1").Clear() | Out-Null

$ws.Range("A1").Value = "CAPID"
$ws.Range("B1").Value = "Squadron Assignment"
$ws.Range("C1").Value = "Flight Assignment"
$ws.Range("D1").Value = "Position"

# Data rows
$ws.Range("A2").Value = 589351
$ws.Range("B2").Value = "Squadron 1"
$ws.Range("C2").Value = "Alpha"
$ws.Range("D2").Value = "Flight Commander"

$ws.Range("A3").Value = 642624
$ws.Range("B3").Value = "Squadron 2"
$ws.Range("C3").Value = "Alpha"
$ws.Range("D3").Value = "Basic Cadet"

$ws.Range("A4").Value = 111111
$ws.Range("B4").Value = "Squadron 3"
$ws.Range("C4").Value = "Bravo"
$ws.Range("D4").Value = "Flight Commander"

$ws.Range("A5").Value = 222222
$ws.Range("B5").Value = "Squadron 4"
$ws.Range("C5").Value = "Bravo"
$ws.Range("D5").Value = "ATS Cadet"

# Column widths: column D inherits the original bestFit width (15.88671875)
# that shifted over from the pre-insert B:C range, so it's left alone. Columns
# B and C get their own explicit widths for the new content (set via the
# closest values this host's width-rounding maps onto the target widths of
# 19.6640625 and 15.88671875 respectively).
$ws.Columns.Item(2).ColumnWidth = 18.8
$ws.Columns.Item(3).ColumnWidth = 15

$ws.Range("C8").Select() | Out-Null
